$wb = $excel.ActiveWorkbook

# Updated "want to go" (想去人数) counts in column F for the sheets that
# carry the event data: "展览" and "全部类型".
$updates = @{
    2  = 6494
    3  = 187
    6  = 1948
    7  = 1486
    8  = 302
    9  = 995
    10 = 357
    12 = 5619
    13 = 74
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
